$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 already carries the plain centered style (s=1) used across the sheet;
# just layer the per-column formatting (number format / wrap) on top, and copy
# the bordered style used in the other rows for column B.

# A9: commit date (15-Jun-2022 = serial 44727), keep centered, no border, custom date format
$ws.Range("A9").Value = 44727
$ws.Range("A9").NumberFormat = "d-mmm"
$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4108

# B9: author, matches the bordered style used by every other row in column B
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Octavio Lucardi Fierro"

# C9: commit title
$ws.Range("C9").Value = "Juan camina bien + menu de opciones"
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").VerticalAlignment = -4108
$ws.Range("C9").WrapText = $true

# D9: commit description
$ws.Range("D9").Value = "Juan tiene bien hechas las animaciones, hay un menu de opciones que regula el sonido y la pantalla completa y se intenta disparar"
$ws.Range("D9").HorizontalAlignment = -4108
$ws.Range("D9").VerticalAlignment = -4108
$ws.Range("D9").WrapText = $true

$ws.Rows.Item(9).RowHeight = 45

$ws.Range("D11").Select()
